# 📊 Horarios actualizados Línea 141 - 812
# Updates the scraped-schedule workbook (3 sheets: LP1912, LP1912-215, 6203-6173)
# with a fresh scrape timestamp (02:21:47) and new/changed rows on the LP1912
# sheet (and the mirrored "combined" row on LP1912-215).

$wb = $excel.ActiveWorkbook

$newTime = "02:21:47"

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A2").Value = "Última actualización: " + $newTime
$ws1.Range("A3").Value = "Total filas: 7"

# Row 8: refresh scrape time + recomputed minutes-until-arrival
$ws1.Range("A8").Value = $newTime
$ws1.Range("D8").Value = 37

# New rows appended to the schedule
$ws1.Range("A11").Value = $newTime
$ws1.Range("B11").Value = "03:56"
$ws1.Range("C11").Value = "14_ABASTO"
$ws1.Range("D11").Value = 95
$ws1.Range("E11").Value = "LP1912"

$ws1.Range("A12").Value = $newTime
$ws1.Range("B12").Value = "04:01"
$ws1.Range("C12").Value = "81_EL PELIGRO"
$ws1.Range("D12").Value = 100
$ws1.Range("E12").Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A2").Value = "Última actualización: " + $newTime

# Row 7: refresh scrape time + recomputed minutes-until-arrival
$ws2.Range("A7").Value = $newTime
$ws2.Range("D7").Value = 37

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A2").Value = "Última actualización: " + $newTime
